# v0.7.3y: GameMap::features is HashMap
# Adds a new "Light Switch" feature (rows 195-198) to the Features sheet,
# tags two existing features with their sub-type in column D, and updates
# the saved cell-selection state on the Terrains / Features sheets.

$wb = $excel.ActiveWorkbook

$features = $wb.Worksheets.Item("Features")
$terrains = $wb.Worksheets.Item("Terrains")

# --- Tag two existing rows with a sub-type in column D -------------------
# Row 68 = feature id 186 "Outside Light Source" -> Light Source
$features.Range("D68").Value2 = "Light Source"
# Row 73 = feature id 191 "Window" -> Furniture
$features.Range("D73").Value2 = "Furniture"

# --- Insert the new "Light Switch" feature (4 rows) before row 77 --------
$features.Rows("77:80").Insert()

$features.Range("B77").Value2 = 195
$features.Range("C77").Value2 = "Light Switch"
$features.Range("D77").Value2 = "Appliance"
$features.Range("E77").Value2 = 0
$features.Range("H77").Value2 = "up"

$features.Range("B78").Value2 = 196
$features.Range("H78").Value2 = "down"

$features.Range("B79").Value2 = 197
$features.Range("H79").Value2 = "left"

$features.Range("B80").Value2 = 198
$features.Range("H80").Value2 = "right"

# --- Restore saved selections (Terrains first, Features last so it stays
#     the active/visible tab, matching the original workbook state) ------
$terrains.Range("C13").Select()
$features.Range("D74").Select()
